{"js": "// Find the paragraph that holds the \"inputDateHeader\" conditional template\n// expression (split across three runs in the original document) and replace\n// its whole text with the merged/updated Jinja-style expression that adds\n// date formatting for both the electronic and paper input-date values.\nconst body = context.document.body;\n\nconst searchText = \"{% if inputDateHeader %}{{ inputDateHeader }}\";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target placeholder text in document body.\");\n}\n\nconst targetRange = results.items[0];\nconst paragraph = targetRange.paragraphs.getFirst();\n\nconst newText =\n  '{% if inputDateHeader %}{{ inputDateHeader | date(\"dd.MM.YYYY\") }}' +\n  '{% if paperInputDateHeader %} ({{ paperInputDateHeader | date(\"dd.MM.YYYY\") }}){% else %}{% endif %}' +\n  '{% else %}-{% endif %}';\n\n// Replace the entire paragraph's text (all 3 runs) with a single run\n// carrying the new merged text, keeping the paragraph's formatting.\nparagraph.insertText(newText, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Locate the paragraph holding the \"inputDateHeader\" conditional template\n# expression (originally split across three runs) and replace its whole\n# text with the merged/updated expression that adds date formatting for\n# both the electronic and paper input-date values.\n\n$d = $word.ActiveDocument\n\n$oldText = \"{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}\"\n$newText = '{% if inputDateHeader %}{{ inputDateHeader | date(\"dd.MM.YYYY\") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date(\"dd.MM.YYYY\") }}){% else %}{% endif %}{% else %}-{% endif %}'\n\n$find = $d.Content.Find\n$find.Text = $oldText\n$found = $find.Execute()\n\nif ($found) {\n    # Assign directly to the matched Range's Text (rather than using\n    # Find's ReplaceWith mechanism) so Word's smart-quote autocorrect\n    # does not mangle the literal straight double quotes in the filter\n    # arguments. This also collapses the paragraph's three runs into a\n    # single run, inheriting the formatting of the first one.\n    $rng = $find.Parent\n    $rng.Text = $newText\n} else {\n    throw \"Could not find target placeholder text in document content.\"\n}\n"}
